$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Global font fix: TimesNewToman -> Times New Roman (affects every run)
# ---------------------------------------------------------------------------
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------------
# Helper: replace all occurrences of $old with $new using a fresh range each
# time (Find leaves its range advanced, so reusing one range across multiple
# calls silently misses earlier matches).
# ---------------------------------------------------------------------------
function Replace-Text($old, $new) {
    $rng = $d.Range(0, $d.Content.End)
    $rng.Find.ClearFormatting()
    [void]$rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# ---------------------------------------------------------------------------
# 2. Title
# ---------------------------------------------------------------------------
Replace-Text "Cultural Mosaic: A Vibrant Tapestry of Human Experience" "Mathematics: The Language of Science and Logic"

# ---------------------------------------------------------------------------
# 3. Author name / handle / e-mail
# ---------------------------------------------------------------------------
Replace-Text "Sarah Morgan" "Albert Spencer"
Replace-Text "sarah" "albert"
Replace-Text "morgan@esquire" "spencer@educonnect"

# The final ".edu" -> ".org" is handled separately (narrow range) because
# "edu" is now also a substring of "educonnect".
$p3 = $d.Paragraphs.Item(3).Range
$innerEnd = $p3.End - 1
$eduStart = $innerEnd - 3
$eduRng = $d.Range($eduStart, $innerEnd)
if ($eduRng.Text -eq "edu") {
    $eduRng.Text = "org"
}

Write-Host "Header block done:" $d.Paragraphs.Item(1).Range.Text "|" $d.Paragraphs.Item(2).Range.Text "|" $d.Paragraphs.Item(3).Range.Text

# ---------------------------------------------------------------------------
# 4. Body paragraph (5th paragraph): full rewrite with new wording and a
#    different line-break layout. We clear the paragraph's content (keeping
#    the paragraph mark) and rebuild it run-by-run / break-by-break so the
#    <w:br/> placement matches the target exactly.
# ---------------------------------------------------------------------------
function Rebuild-Paragraph($paraIndex, $items) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $start = $r.Start
    $inner = $d.Range($start, $r.End - 1)
    $inner.Text = ""

    $pos = $start
    foreach ($item in $items) {
        if ($item.br) {
            $insBr = $d.Range($pos, $pos)
            $insBr.InsertBreak(6)
            $pos = $pos + 1
        }
        if ($item.text.Length -gt 0) {
            $insText = $d.Range($pos, $pos)
            $insText.InsertAfter($item.text)
            $pos = $pos + $item.text.Length
        }
    }

    $full = $d.Range($start, $pos)
    $full.Font.Name = "Times New Roman"
    $full.Font.Color = 0
    $full.Font.Size = 12
}

$body1Items = @(
    @{ br = $false; text = "Mathematics, the intricate language of numbers and logic, pervades our universe like an omnipresent symphony" },
    @{ br = $false; text = "." },
    @{ br = $false; text = " From the mesmerizing patterns of nature to the profound workings of the cosmos, mathematics weaves an invisible tapestry of order and harmony" },
    @{ br = $false; text = "." },
    @{ br = $false; text = " Unraveling its enigmatic secrets, we embark on an extraordinary quest for knowledge and understanding, shedding light on the universe's hidden symmetries and revealing the elegance of its underlying principles" },
    @{ br = $false; text = "." },
    @{ br = $true;  text = "" },
    @{ br = $true;  text = "In the realm of science, mathematics reigns supreme as an indispensable tool for unraveling the mysteries of nature" },
    @{ br = $false; text = "." },
    @{ br = $false; text = " It provides the language and framework through which we quantify and analyze data, formulate hypotheses, and construct theories" },
    @{ br = $false; text = "." },
    @{ br = $false; text = " From the complex equations governing the motion of celestial bodies to the intricate interactions of subatomic particles, mathematics empowers us to comprehend the intricate workings of the universe, uncovering its fundamental laws and unlocking its profound secrets" },
    @{ br = $false; text = "." },
    @{ br = $true;  text = "" },
    @{ br = $true;  text = "Beyond its scientific applications, mathematics also plays a vital role in our everyday lives" },
    @{ br = $false; text = "." },
    @{ br = $false; text = " From calculating percentages in financial transactions to deciphering complex graphs and charts in news articles, mathematical literacy is an essential skill for informed decision-making and critical thinking" },
    @{ br = $false; text = "." },
    @{ br = $false; text = " Whether navigating intricate traffic patterns or managing personal budgets, mathematics equips us with the tools to navigate the complexities of modern life, empowering us to make informed choices and solve real-world problems" },
    @{ br = $false; text = "." }
)

Rebuild-Paragraph 5 $body1Items

Write-Host "Paragraph 5 rebuilt, length:" $d.Paragraphs.Item(5).Range.Text.Length

# ---------------------------------------------------------------------------
# 5. Summary paragraph (7th paragraph): sentence-for-sentence rewrite, no
#    change to the number/placement of line breaks, so plain Find/Replace
#    (scoped to this paragraph) is sufficient and keeps things simple.
# ---------------------------------------------------------------------------
function Replace-InParagraph($paraIndex, $old, $new) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $d.Range($p.Range.Start, $p.Range.End)
    [void]$rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-InParagraph 7 "Cultural diversity, akin to a vibrant mosaic, constitutes the multifaceted story of human civilization" "Mathematics, the language of science and logic, is a gateway to understanding the universe's profound secrets"
Replace-InParagraph 7 " It influences our way of life, from our beliefs and values to our customs and traditions" " It enables us to quantify and analyze data, formulate hypotheses, and construct theories, providing a framework for comprehending the intricate workings of nature"
Replace-InParagraph 7 " The rich tapestry of cultures propels innovation and creativity while fostering understanding and empathy among individuals and communities" " Beyond its scientific applications, mathematics empowers us to make informed decisions, solve real-world problems, and navigate the complexities of modern life"
Replace-InParagraph 7 " As citizens of a globalized world, it is imperative to not only appreciate cultural diversity but also actively contribute to its preservation and promotion" " As we delve deeper into the enigmatic world of mathematics, we unlock the keys to unlocking the universe's mysteries and harnessing its knowledge for the betterment of humanity"

Write-Host "Summary:" $d.Paragraphs.Item(7).Range.Text

# ---------------------------------------------------------------------------
# 6. A new, empty trailing paragraph is added right before the section break.
# ---------------------------------------------------------------------------
$endRng = $d.Range($d.Content.End, $d.Content.End)
$endRng.InsertParagraphAfter()

Write-Host "Final paragraph count:" $d.Paragraphs.Count


